$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column (02-dec) before column DY (129) ---
$prixSpot = $wb.Worksheets.Item(1)
$prixSpot.Columns.Item(129).Insert()

# New header for the inserted column
$prixSpot.Cells.Item(1, 129).Value = "02-dec"

# New column has no data yet for this day (same placeholder used elsewhere in the sheet)
for ($r = 2; $r -le 25; $r++) {
    $prixSpot.Cells.Item($r, 129).Value = "-"
}

# --- Sheet "Gaz": append two new daily rows ---
$gaz = $wb.Worksheets.Item(2)
$gaz.Cells.Item(156, 1).Formula = "'2025-11-29"
$gaz.Cells.Item(156, 2).Value = 27.525
$gaz.Cells.Item(157, 1).Formula = "'2025-11-30"
$gaz.Cells.Item(157, 2).Value = 27.525

# Reuse the formatting (plain text / plain number, no quote-prefix flag) from the row above
$gaz.Cells.Item(155, 1).Copy()
$gaz.Range("A156:A157").PasteSpecial(-4122)
$gaz.Cells.Item(155, 2).Copy()
$gaz.Range("B156:B157").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Sheet "CO2": append two new daily rows (values not yet available) ---
$co2 = $wb.Worksheets.Item(3)
$co2.Cells.Item(157, 1).Formula = "'2025-11-29"
$co2.Cells.Item(158, 1).Formula = "'2025-11-30"

$co2.Cells.Item(156, 1).Copy()
$co2.Range("A157:A158").PasteSpecial(-4122)
$excel.CutCopyMode = 0
